$d = $word.ActiveDocument

function Do-Break($findText, $replaceText) {
    $result = $d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $result) {
        throw "Find/Replace failed for: $findText"
    }
}

# --- Programa (PT) paragraph: insert line breaks between topics ---
Do-Break "(4 horas)ESTEQUIOMETRIA" "(4 horas)^lESTEQUIOMETRIA"
Do-Break "(4 horas)REAÇÕES A VOLUME CONSTANTE" "(4 horas)^lREAÇÕES A VOLUME CONSTANTE"
Do-Break "(12 horas)REAÇÕES A VOLUME VARIÁVEL" "(12 horas)^lREAÇÕES A VOLUME VARIÁVEL"
Do-Break "(4 horas)INTRODUÇÃO A PROJETO DE REATORES" "(4 horas)^lINTRODUÇÃO A PROJETO DE REATORES"
Do-Break "(8 horas)COLETA E ANÁLISE DE DADOS CINÉTICOS:" "(8 horas)^lCOLETA E ANÁLISE DE DADOS CINÉTICOS:"
Do-Break "CINÉTICOS:Métodos diferencial" "CINÉTICOS:^lMétodos diferencial"
Do-Break "(12 horas)CINÉTICA DAS REAÇÕES COMPLEXAS" "(12 horas)^lCINÉTICA DAS REAÇÕES COMPLEXAS"
Do-Break "(8 horas)CATÁLISE" "(8 horas)^lCATÁLISE"

# --- Programa (EN, italic) paragraph: insert line breaks between topics ---
Do-Break "(4 hours)KINETIC STOICHIOMETRY" "(4 hours)^lKINETIC STOICHIOMETRY"
Do-Break "(4 hours)REACTIONS AT CONSTANT VOLUME" "(4 hours)^lREACTIONS AT CONSTANT VOLUME"
Do-Break "(12 hours)REACTIONS AT VARIABLE VOLUME" "(12 hours)^lREACTIONS AT VARIABLE VOLUME"
Do-Break "(4 hours)INTRODUCTION TO REACTOR DESIGN" "(4 hours)^lINTRODUCTION TO REACTOR DESIGN"
Do-Break "(PFR). (8 hours)COLLECTION AND ANALYSIS OF KINETIC DATA" "(PFR). (8 hours)^lCOLLECTION AND ANALYSIS OF KINETIC DATA"
Do-Break "(12 hours)KINETICS OF COMPLEX REACTIONS" "(12 hours)^lKINETICS OF COMPLEX REACTIONS"
Do-Break "(8 hours)CATALYSIS" "(8 hours)^lCATALYSIS"

# --- Critério run: split sentence with a line break ---
Do-Break "50% P2.Obs: fica a critério" "50% P2.^lObs: fica a critério"

# --- Bibliografia paragraph: insert line breaks between references ---
Do-Break "Bibliografia Básica:1- FOGLER" "Bibliografia Básica:^l1- FOGLER"
Do-Break "2009.2- LEVENSPIEL" "2009.^l2- LEVENSPIEL"
Do-Break "2000)3- VAN SANTEN" "2000)^l3- VAN SANTEN"
Do-Break "1995.Bibliografia Complementar:1- Missen" "1995.^l^lBibliografia Complementar:^l1- Missen"
Do-Break "1999.2- Rothenberg" "1999.^l2- Rothenberg"
Do-Break "Chichester.3- DENISOV" "Chichester.^l3- DENISOV"
Do-Break "2003.4- Hagen" "2003.^l4- Hagen"
Do-Break "2006.5- Salmi" "2006.^l5- Salmi"
Do-Break "2011.6- Mortimer" "2011.^l6- Mortimer"
Do-Break "2002.7- FROMENT" "2002.^l7- FROMENT"
Do-Break "1990.8- HILL" "1990.^l8- HILL"
Do-Break "1977.9- SMITH" "1977.^l9- SMITH"
Do-Break "1981.10- DENBIGH" "1981.^l10- DENBIGH"
Do-Break "1970.11 - AGUIAR" "1970.^l11 - AGUIAR"

Write-Host "Done."
